$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average_county_temperature (column I) values using NOAA data
$ws.Range("I7").Value = -0.763888888888889
$ws.Range("I12").Value = 12.41429539295394
$ws.Range("I13").Value = 13.75752314814816
$ws.Range("I15").Value = 21.79166666666666
$ws.Range("I16").Value = 13.75752314814816
$ws.Range("I23").Value = 20.68981481481483
$ws.Range("I26").Value = -0.763888888888889
$ws.Range("I27").Value = 12.67039049919483
$ws.Range("I28").Value = 19.65277777777778
